# Updated symbol list (price + 1h volume change) per commit diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "258.51"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "0.81%"

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "26.86"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "-0.71%"

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "4.638"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "-0.29%"

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.05956"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "1.11%"

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "6.632"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "-0.03%"

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.8564"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "-0.77%"

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.9238"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "-0.76%"

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.1387"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "-1.37%"

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.04262"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "14.94%"

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07019"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "-1.03%"

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.02961"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "-8.25%"

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.09106"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "-1.13%"

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.001540"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "0.37%"

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.0006050"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "0.35%"

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.006121"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "0.65%"

$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "-1.79%"

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "3.125"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "-2.03%"

$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "-2.17%"

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.3105"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "0.21%"

$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "1.66%"

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "3.933"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "2.23%"

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.04224"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "-0.20%"

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.001217"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "-0.33%"

$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "4.68%"

$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "-0.09%"

$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "-11.60%"

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.03824"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "0.01%"

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.1109"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "0.80%"

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.003784"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "-4.14%"

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.002428"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "0.74%"

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.01494"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "31.25%"

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.00005152"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "-5.86%"

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.00000000749"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "-0.06%"

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.04996"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "-17.00%"

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.2380"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "10,344.99%"

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.00002098"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "-0.06%"

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.0001998"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "-0.06%"
